$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for columns D, M, N, P, S
$ws.Range("D2").Value = 44991
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 6000
$ws.Range("P2").Value = 6000
$ws.Range("S2").Value = 3000

$ws.Range("D3").Value = 44995
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 5500
$ws.Range("P3").Value = 5750
$ws.Range("S3").Value = 2875
